# Clean up the Sheet1 header row: the original headers carried extraction
# annotations ("BREWERY[0](STRING)", "WEBSITE[1] (STRING)", etc.) left over
# from the scraping/indexing step. Replace them with plain column titles.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "BREWERY"
$ws.Range("C1").Value = "WEBSITE"
$ws.Range("D1").Value = "STREET ADDRESS"
$ws.Range("E1").Value = "CITY"
$ws.Range("F1").Value = "STATE"
$ws.Range("G1").Value = "ZIPCODE"
$ws.Range("H1").Value = "PHONE"
$ws.Range("I1").Value = "DAYS CLOSED"

# Leave the selection on I1, matching the saved cursor position.
$ws.Range("I1").Select() | Out-Null
